$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the R10 rule row from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Move/restore the active selection to E8, matching the saved view state
$ws.Range("E8").Select()
